$wb = $excel.ActiveWorkbook

# --- Update the config sheet's referenced template paths ---
$ws = $wb.Worksheets.Item("config")
$ws.Activate()

$ws.Range("B7").Value = "../5MW_Baseline/wind/TurbSim.inp"
$ws.Range("B8").Value = "../5MW_Baseline/wind/IEC_template.IPT"

# Reflect the new selection (user's cursor moved onto the cells just edited)
$ws.Range("B7:B8").Select()
